$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 15, shifting existing rows 15-22 down to 16-23 ---
$ws.Range("A15").EntireRow.Insert()

# --- Populate new row 15 (Clear Lake 23773411 2,683,868 m3/day) ---
$ws.Range("A15").Value = "CW3M"
$ws.Range("B15").Value = "Baseline_2010-18_C91"
$ws.Range("C15").Value = 2010
$ws.Range("D15").Value = 1090.199341000000003987224772572517
$ws.Range("E15").Value = 1990.467650999999932537321001291275
$ws.Range("F15").Value = 1.255063000000000039690917219559
$ws.Range("G15").Value = 347.021851000000026488123694434762
$ws.Range("H15").Value = 10.610913999999999290935193130281
$ws.Range("I15").Value = 3.822232000000000073924866228481
$ws.Range("J15").Value = 8.840457000000000675754563417286
$ws.Range("K15").Value = 814.495055999999976847902871668339
$ws.Range("L15").Value = 93.229797000000004913999873679131
$ws.Range("M15").Value = 1324.581421000000091225956566631794
$ws.Range("N15").Value = 1201.767211999999972249497659504414
$ws.Range("O15").Value = 7126.601563000000169267877936363220
$ws.Range("P15").Value = 29450.638672000000951811671257019043
$ws.Range("Q15").Value = -0.463108000000000019635848502730
$ws.Range("R15").Value = -0.000134000000000000004524158825
$ws.Range("S15").Value = 2010

# --- Append new row 24 (Blue R at and above Tidbits Cr 23773429 138,240 m3/day) ---
$ws.Range("A24").Value = "CW3M"
$ws.Range("B24").Value = "Baseline_2010-18_C91"
$ws.Range("C24").Value = "2010-18"
$ws.Range("D24:N24").NumberFormat = "0.00"
$ws.Range("Q24").NumberFormat = "0.00"
$ws.Range("O24:P24").NumberFormat = "0"
$ws.Range("R24").NumberFormat = "0.000000"
$ws.Range("D24").Value = 1161.159905444444348177057690918446
$ws.Range("E24").Value = 1901.515733444444322230992838740349
$ws.Range("F24").Value = 1.011925555555555744646767379891
$ws.Range("G24").Value = 347.233127444444448883587028831244
$ws.Range("H24").Value = 9.775355222222223972039500949904
$ws.Range("I24").Value = 4.681360777777778814368048188044
$ws.Range("J24").Value = 8.145128999999998953285285097081
$ws.Range("K24").Value = 769.369927111111110207275487482548
$ws.Range("L24").Value = 83.470620444444449503862415440381
$ws.Range("M24").Value = 1403.060424888888974237488582730293
$ws.Range("N24").Value = 1161.285651444444283697521314024925
$ws.Range("O24").Value = 4883.927707333332364214584231376648
$ws.Range("P24").Value = 27227.338324888889474095776677131653
$ws.Range("Q24").Value = -0.045654777777777787173274504084
$ws.Range("R24").Value = -0.000037111111111111106655377490
$ws.Range("S24").Value = "2010-18"

# --- Update selection to match saved view state ---
$ws.Range("A15:B15").Select()

